# Add a new worksheet "Test Cases of Add to Cart" between the existing
# "Test Cases for Login " sheet and the "Bug Report for Login" sheet,
# populate it with three new manual test cases, and refresh a few
# cosmetic details (tab colors + selections) on all three sheets.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item(1)
$bugSheet   = $wb.Worksheets.Item("Bug Report for Login")

# --- "Test Cases for Login " sheet: move the selection, no longer scrolled ---
$loginSheet.Range("D2").Select()

# Tab colors (VBA-style packed BGR long = R + G*256 + B*65536)
$loginSheet.Tab.Color = 49407    # FFC000 (orange)
$bugSheet.Tab.Color   = 5296274  # 92D050 (green)

# --- Insert the new sheet right before "Bug Report for Login" ---
$bugSheet.Activate()
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Test Cases of Add to Cart"
$newSheet.Tab.Color = 49407      # FFC000 (orange), same as the login sheet

# Copy the header/body cell formatting from the login sheet so the new
# sheet reuses the existing styles instead of creating new ones.
$loginSheet.Range("A1:G1").Copy()
$newSheet.Range("A1:G1").PasteSpecial(-4122)   # xlPasteFormats

$loginSheet.Range("A2:G2").Copy()
$newSheet.Range("A2:G4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Header row ---
$newSheet.Range("A1").Value = "Test Case ID"
$newSheet.Range("B1").Value = "Test Scenario"
$newSheet.Range("C1").Value = "Precondition"
$newSheet.Range("D1").Value = "Steps"
$newSheet.Range("E1").Value = "Expected Result"
$newSheet.Range("F1").Value = "Priority"
$newSheet.Range("G1").Value = "Status"

# --- TC009 ---
$newSheet.Range("A2").Value = "TC009"
$newSheet.Range("B2").Value = "Add single to cart "
$newSheet.Range("C2").Value = "User logged in (standard_user)on inventory page"
$newSheet.Range("D2").Value = "1-Open inventory`n2-Click Add to cart on first product"
$newSheet.Range("E2").Value = "Product is added ,cart badge shows (1)product appear in cart page"
$newSheet.Range("F2").Value = "High"
$newSheet.Range("G2").Value = "Pass"

# --- TC010 ---
$newSheet.Range("A3").Value = "TC010"
$newSheet.Range("B3").Value = "Add all items to cart "
$newSheet.Range("C3").Value = "Userlogged in on inventory page "
$newSheet.Range("D3").Value = "1-Click Add to cart for every product"
$newSheet.Range("E3").Value = "Cart badge number equals total products(6);all items listed in Cart"
$newSheet.Range("F3").Value = "High"
$newSheet.Range("G3").Value = "Pass"

# --- TC011 ---
$newSheet.Range("A4").Value = "TC011"
$newSheet.Range("B4").Value = "Remove item from cart"
$newSheet.Range("C4").Value = "At least one item in cart "
$newSheet.Range("D4").Value = "1-Go to Cart`n2-Click Remove on a product"
$newSheet.Range("E4").Value = "Item removed,cart badge decreased product not listed in cart "
$newSheet.Range("F4").Value = "High"
$newSheet.Range("G4").Value = "Pass"

# --- Row heights to match the source layout ---
$newSheet.Rows.Item(1).RowHeight = 56.25
$newSheet.Rows.Item(2).RowHeight = 30
$newSheet.Rows.Item(3).RowHeight = 26.25
$newSheet.Rows.Item(4).RowHeight = 30

# --- Column widths to match the source layout ---
$newSheet.Columns.Item(1).ColumnWidth = 22.140625
$newSheet.Columns.Item(2).ColumnWidth = 21.7109375
$newSheet.Columns.Item(3).ColumnWidth = 45.140625
$newSheet.Columns.Item(4).ColumnWidth = 33.28515625
$newSheet.Columns.Item(5).ColumnWidth = 60.7109375
$newSheet.Columns.Item(6).ColumnWidth = 9.85546875
$newSheet.Columns.Item(7).ColumnWidth = 10.7109375

# Final selection on the new sheet
$newSheet.Range("A9").Select()

Write-Output "Workbook now has $($wb.Worksheets.Count) sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output " - $($s.Name)"
}
